# "Revised data files to make them R-friendly"
#
# The original workbook had headers "Female"/"Male" in A1/B1 and the
# worksheet still carried an interactive UI selection (cell G10) left
# over from whoever last edited the file in Excel. Re-saving the file
# through a script (e.g. R's readxl/openxlsx workflow) lower-cases the
# header labels and drops the stale UI selection, leaving the view
# pointed at the top-left (default) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lower-case the category labels used as column headers.
$ws.Range("A1").Value = "female"
$ws.Range("B1").Value = "male"

# Clear the leftover "G10" selection/scroll position that was saved in
# the original file; move the view back to the default A1 cell.
$ws.Range("A1").Select() | Out-Null
